$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for all data rows (2..138)
# from the old serial date 45172 (2023-09-03) to the new serial date 45175 (2023-09-06).
$ws.Range("C2:C138").Value = 45175
